$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 0.90239999999999998
$ws.Range("E3").Value = 0.90669999999999995
$ws.Range("G3").Value = 0.90649999999999997
$ws.Range("C4").Value = 0.98680000000000001
$ws.Range("D4").Value = 0.98629999999999995
$ws.Range("E4").Value = 0.98729999999999996
$ws.Range("F4").Value = 0.9859
$ws.Range("G4").Value = 0.98729999999999996
$ws.Range("H4").Value = 0.98670000000000002
$ws.Range("D5").Value = 9.1300000000000008
$ws.Range("F5").Value = 4.41
$ws.Range("H5").Value = 4.59
$ws.Range("D6").Value = 44.51
$ws.Range("F6").Value = 22.14
$ws.Range("G6").Value = 0.2
$ws.Range("H6").Value = 22.4
$ws.Range("D7").Value = -1.76
$ws.Range("F7").Value = -0.67
$ws.Range("H7").Value = -0.78
$ws.Range("C9").Value = 0.88670000000000004
$ws.Range("D9").Value = 0.92090000000000005
$ws.Range("E9").Value = 0.89180000000000004
$ws.Range("F9").Value = 0.91920000000000002
$ws.Range("G9").Value = 0.89149999999999996
$ws.Range("H9").Value = 0.92159999999999997
$ws.Range("D10").Value = 23.4
$ws.Range("F10").Value = 11.43
$ws.Range("G10").Value = 0.32
$ws.Range("H10").Value = 11.91
$ws.Range("C11").Value = 0.7
$ws.Range("D11").Value = 106.76
$ws.Range("F11").Value = 52.8
$ws.Range("G11").Value = 0.6
$ws.Range("H11").Value = 54.44
$ws.Range("D12").Value = -1.96
$ws.Range("F12").Value = -0.59
$ws.Range("H12").Value = -0.72
$ws.Range("E13").Value = 0.86
$ws.Range("G13").Value = 0.86
$ws.Range("C14").Value = 0.76749999999999996
$ws.Range("D14").Value = 0.82240000000000002
$ws.Range("E14").Value = 0.77829999999999999
$ws.Range("F14").Value = 0.81840000000000002
$ws.Range("G14").Value = 0.77649999999999997
$ws.Range("H14").Value = 0.81940000000000002
$ws.Range("C15").Value = 0.61
$ws.Range("D15").Value = 35.81
$ws.Range("E15").Value = 0.51
$ws.Range("F15").Value = 17.61
$ws.Range("H15").Value = 18.559999999999999
$ws.Range("C16").Value = 1.01
$ws.Range("D16").Value = 156
$ws.Range("F16").Value = 77.760000000000005
$ws.Range("G16").Value = 0.87
$ws.Range("H16").Value = 80.16
$ws.Range("C17").Value = 0.42
$ws.Range("D17").Value = 2.5099999999999998
$ws.Range("E17").Value = 0.34
$ws.Range("F17").Value = 1.97
$ws.Range("H17").Value = 1.83
$ws.Range("D18").Value = 0.87
$ws.Range("F18").Value = 0.88
$ws.Range("H18").Value = 0.88
$ws.Range("C19").Value = 0.91010000000000002
$ws.Range("D19").Value = 0.89970000000000006
$ws.Range("E19").Value = 0.91449999999999998
$ws.Range("F19").Value = 0.89610000000000001
$ws.Range("G19").Value = 0.91420000000000001
$ws.Range("H19").Value = 0.90249999999999997
$ws.Range("C20").Value = 0.35
$ws.Range("D20").Value = 24.29
$ws.Range("E20").Value = 0.28999999999999998
$ws.Range("F20").Value = 11.83
$ws.Range("H20").Value = 12.29
$ws.Range("C21").Value = 0.62
$ws.Range("D21").Value = 119.43
$ws.Range("F21").Value = 59.63
$ws.Range("H21").Value = 60.22
$ws.Range("D22").Value = -1.39
$ws.Range("E22").Value = 0.12
$ws.Range("F22").Value = -0.06
$ws.Range("G22").Value = 0.12
$ws.Range("H22").Value = -0.38
$ws.Range("E23").Value = 0.86
$ws.Range("C24").Value = 0.93610000000000004
$ws.Range("D24").Value = 0.95279999999999998
$ws.Range("E24").Value = 0.93810000000000004
$ws.Range("F24").Value = 0.95079999999999998
$ws.Range("G24").Value = 0.93989999999999996
$ws.Range("H24").Value = 0.95320000000000005
$ws.Range("C25").Value = 0.28999999999999998
$ws.Range("D25").Value = 17.559999999999999
$ws.Range("F25").Value = 8.68
$ws.Range("H25").Value = 8.94
$ws.Range("C26").Value = 0.52
$ws.Range("D26").Value = 83.16
$ws.Range("E26").Value = 0.45
$ws.Range("F26").Value = 41.54
$ws.Range("G26").Value = 0.44
$ws.Range("H26").Value = 42.42
$ws.Range("D27").Value = -2.1
$ws.Range("F27").Value = -0.73
$ws.Range("H27").Value = -0.82
$ws.Range("C28").Value = 0.89
$ws.Range("H28").Value = 0.88
$ws.Range("C29").Value = 0.92920000000000003
$ws.Range("D29").Value = 0.84870000000000001
$ws.Range("E29").Value = 0.93100000000000005
$ws.Range("F29").Value = 0.83320000000000005
$ws.Range("G29").Value = 0.93310000000000004
$ws.Range("H29").Value = 0.86109999999999998
$ws.Range("C30").Value = 0.31
$ws.Range("D30").Value = 23.21
$ws.Range("E30").Value = 0.26
$ws.Range("F30").Value = 11.38
$ws.Range("G30").Value = 0.26
$ws.Range("H30").Value = 11.59
$ws.Range("C31").Value = 0.54
$ws.Range("D31").Value = 151.69
$ws.Range("F31").Value = 77.010000000000005
$ws.Range("G31").Value = 0.46
$ws.Range("H31").Value = 74.66
$ws.Range("C32").Value = 0.18
$ws.Range("D32").Value = 6.03
$ws.Range("F32").Value = 3.61
$ws.Range("G32").Value = 0.15
$ws.Range("H32").Value = 3.25
$ws.Range("C33").Value = 0.88
$ws.Range("D33").Value = 0.7
$ws.Range("G33").Value = 0.89
$ws.Range("H33").Value = 0.72

$ws.Activate()
$ws.Range("G11").Select()
